# Journal de bord - "JB & Uses Cases & Scénarios" update
# Fills in rows 13-20 of the log table (dates, tasks, durations, remarks)
# that were previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 13, 14 and 15 were still using the "early" cell style
#    (s=1/3/2); bring them in line with the rest of the filled table
#    (s=8/9/10) by copying the formatting from an already-styled row.
#    Row 16 already has the s=1/3/2 style in the source file, so it is
#    left untouched.
# ------------------------------------------------------------------
$ws.Range("A17:D17").Copy($ws.Range("A13:D13"))
$ws.Range("A17:D17").Copy($ws.Range("A14:D14"))
$ws.Range("A17:D17").Copy($ws.Range("A15:D15"))

# ------------------------------------------------------------------
# 2) Fill in the values. The order below reproduces the order in
#    which the rows were actually populated by the author (row 14
#    downward first, row 13 -- "ABSENT" -- filled in afterwards).
# ------------------------------------------------------------------

# Row 14 - 08/02/2018 - Finalisation du Trello
$ws.Range("A14").Value = 43139
$ws.Range("B14").Value = "Finalisation du Trello"
$ws.Range("C14").Value = 0.020833333333333332
$ws.Range("D14").Value = "Ajout des différents sprint, description etc.."

# Row 15 - 08/02/2018 - Ajout des documents dans Github
$ws.Range("A15").Value = 43139
$ws.Range("B15").Value = "Ajout des documents dans Github"
$ws.Range("C15").Value = 0.003472222222222222

# Row 16 - 08/02/2018 - Création de la maquette
$ws.Range("A16").Value = 43139
$ws.Range("B16").Value = "Création de la maquette + comparaison avec d'autres application"
$ws.Range("C16").Value = 0.09375

# Row 13 - 06/02/2018 - ABSENT
$ws.Range("A13").Value = 43137
$ws.Range("B13").Value = "ABSENT"
$ws.Range("C13").Value = 0.125

# Row 17 - 08/02/2018 - Premier jet des maquettes
$ws.Range("A17").Value = 43139
$ws.Range("B17").Value = "Premier jet des maquettes"
$ws.Range("C17").Value = 0.03125

# Row 18 - 08/02/2018 - Création des uses cases & scénarios
$ws.Range("A18").Value = 43139
$ws.Range("B18").Value = "Création des uses cases & scénarios"
$ws.Range("C18").Value = 0.0625

# Row 19 - 08/02/2018 - Finition des uses cases & scénarios
$ws.Range("A19").Value = 43139
$ws.Range("B19").Value = "Finition des uses cases & scénarios"
$ws.Range("C19").Value = 0.020833333333333332

# Row 20 - 08/02/2018 - Création d'un Hello World sous Android Studio
$ws.Range("A20").Value = 43139
$ws.Range("B20").Value = "Création d'un Hello World sous Android Studio"
$ws.Range("C20").Value = 0.03125

# ------------------------------------------------------------------
# 3) Rows 14 and 16 wrap onto two lines in the "Remarque"/"Tâches
#    effectuée" columns, so Excel grew their height to 30pt.
# ------------------------------------------------------------------
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

# ------------------------------------------------------------------
# 4) Leave the selection where the author ended up (B21) before
#    saving the workbook.
# ------------------------------------------------------------------
$ws.Range("B21").Select()

# The "Total" formula in C35 (=SUM(C2:C34)) recalculates automatically
# to include the newly entered durations.
